$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)   # workflow_config
$ws3 = $wb.Worksheets.Item(3)   # analysis_config

# --- workflow_config: insert a new "track_region" parameter row at row 12 ---
# Insert a blank row, pushing existing rows (ensembl_version etc.) down by one.
$ws2.Rows.Item(12).Insert()

# Copy the formatting (border/fill/font/number-format) from the row that is
# now directly below (the old row 12, now row 13) onto the new blank row so
# the new row matches the sheet's normal parameter-row styling.
$ws2.Range("A13:C13").Copy()
$ws2.Range("A12:C12").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Rows.Item(12).RowHeight = 32

# Fill in the new row's content.
$ws2.Range("A12").Value = "track_region"
$ws2.Range("C12").Value = "chr1:750000-1050000"

$descCell = $ws2.Range("B12")
$descCell.Value = "Genomic region to plot genome tracks over. Example: chr1:500000-900000 (REQUIRED if run_genome_track)"
$descCell.Characters(44, 27).Font.Bold = $true
$descCell.Characters(72, 30).Font.Bold = $true

# --- view/selection bookkeeping to mirror the author's interactive edit ---
$ws3.Activate()
$ws3.Range("C3").Select()

$ws2.Activate()
$ws2.Range("C13").Select()
